$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "23.831.52"
$c.ClearFormats()
$ws.Range("E2").Value = "  -2.88%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.616.66"
$c.ClearFormats()
$ws.Range("E3").Value = "  -3.23%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$ws.Range("E4").Value = "  -0.23%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "307.77"
$c.ClearFormats()
$ws.Range("E5").Value = "  -1.96%  "
$ws.Range("E6").Value = "  -0.23%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.3927"
$c.ClearFormats()
$ws.Range("E7").Value = "  -0.66%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3834"
$c.ClearFormats()
$ws.Range("E8").Value = "  -2.73%  "
$ws.Range("E9").Value = "  -0.35%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "1.355"
$c.ClearFormats()
$ws.Range("E10").Value = "  -2.63%  "
$ws.Range("E11").Value = "  -1.94%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.08442"
$c.ClearFormats()
$ws.Range("E12").Value = "  -2.27%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "23.65"
$c.ClearFormats()
$ws.Range("E13").Value = "  -6.94%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "7.040"
$c.ClearFormats()
$ws.Range("E14").Value = "  -3.67%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "7.557"
$c.ClearFormats()
$ws.Range("E15").Value = "  -1.68%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.00001277"
$c.ClearFormats()
$ws.Range("E16").Value = "  -2.86%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "1.619.04"
$c.ClearFormats()
$ws.Range("E17").Value = "  -3.52%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "93.75"
$c.ClearFormats()
$ws.Range("E18").Value = "  -0.09%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.06930"
$c.ClearFormats()
$ws.Range("E19").Value = "  -1.19%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "19.96"
$c.ClearFormats()
$ws.Range("E20").Value = "  -5.55%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "6.810"
$c.ClearFormats()
$ws.Range("E21").Value = "  -3.73%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "1.001"
$c.ClearFormats()
$ws.Range("E22").Value = "  -0.13%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "13.41"
$c.ClearFormats()
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "23.837.58"
$c.ClearFormats()
$ws.Range("E24").Value = "  -2.85%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.452"
$c.ClearFormats()
$ws.Range("E25").Value = "  +4.60%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "2.835"
$c.ClearFormats()
$ws.Range("E26").Value = "  +2.88%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "22.19"
$c.ClearFormats()
$ws.Range("E27").Value = "  -3.51%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "157.00"
$c.ClearFormats()
$ws.Range("E28").Value = "  -1.82%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "139.80"
$c.ClearFormats()
$ws.Range("E29").Value = "  -4.07%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "5.277"
$c.ClearFormats()
$ws.Range("E30").Value = "  -9.68%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "7.776"
$c.ClearFormats()
$ws.Range("E31").Value = "  -6.16%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "2.483"
$c.ClearFormats()
$ws.Range("E32").Value = "  -1.86%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "1.799.76"
$c.ClearFormats()
$ws.Range("E33").Value = "  -3.26%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.08080"
$c.ClearFormats()
$ws.Range("E34").Value = "  -1.97%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.9725"
$c.ClearFormats()
$ws.Range("E35").Value = "  -1.80%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02887"
$c.ClearFormats()
$ws.Range("E36").Value = "  -6.30%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.559"
$c.ClearFormats()
$ws.Range("E37").Value = "  -5.01%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.2662"
$c.ClearFormats()
$ws.Range("E38").Value = "  -5.06%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.09121"
$c.ClearFormats()
$ws.Range("E39").Value = "  -5.32%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "10.31"
$c.ClearFormats()
$ws.Range("E40").Value = "  +0.30%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "13.59"
$c.ClearFormats()
$ws.Range("E41").Value = "  +0.30%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "1.428"
$c.ClearFormats()
$ws.Range("E42").Value = "  -6.12%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.7482"
$c.ClearFormats()
$ws.Range("E43").Value = "  -4.92%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "16.06"
$c.ClearFormats()
$ws.Range("E44").Value = "  -3.22%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.6903"
$c.ClearFormats()
$ws.Range("E45").Value = "  -2.60%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.462"
$c.ClearFormats()
$ws.Range("E46").Value = "  -3.75%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "4.070"
$c.ClearFormats()
$ws.Range("E47").Value = "  -2.52%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.9999"
$c.ClearFormats()
$ws.Range("E48").Value = "  -0.17%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.08236"
$c.ClearFormats()
$ws.Range("E49").Value = "  -4.54%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "134.65"
$c.ClearFormats()
$ws.Range("E50").Value = "  -2.48%  "
$ws.Range("E51").Value = "  -8.94%  "
